$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78 (existing rows 78-114 shift down to 79-115).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new price record.
$ws.Cells.Item(78, 1).Value  = 7
$ws.Cells.Item(78, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(78, 3).Value  = "Ñuble"
$ws.Cells.Item(78, 4).Value  = 45229
$ws.Cells.Item(78, 5).Value  = 16
$ws.Cells.Item(78, 6).Value  = 100112022
$ws.Cells.Item(78, 7).Value  = "Arveja Verde"
$ws.Cells.Item(78, 8).Value  = "Sin especificar"
$ws.Cells.Item(78, 9).Value  = "Primera"
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 26000
$ws.Cells.Item(78, 12).Value = 26000
$ws.Cells.Item(78, 13).Value = 26000
$ws.Cells.Item(78, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(78, 16).Value = 1040
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"

# Match the date-number-format style used by column D in the other rows.
$ws.Cells.Item(78, 4).NumberFormat = $ws.Cells.Item(79, 4).NumberFormat
